# "Added last minute updates"
#
# The first paragraph of the document holds a merge-field style
# placeholder (docx4j/OpenDoPE bookmark text) that needs to be
# retargeted from the generic "pgi_5301_topic_35" id to the specific
# "AFMC_PGI_5301_9001_93" id, and the paragraph needs to pick up the
# same box-spacing/indent formatting used by the rest of the document's
# body paragraphs. The old trailing " " run is dropped entirely.

$d = $word.ActiveDocument

# Replace "**ID__AFFARS_pgi_5301_topic_35__ID** " (text + the trailing
# space that lived in its own run) with the new placeholder text (no
# trailing space) in a single Find/Replace so Word collapses the
# paragraph back down to one run.
$d.Content.Find.Execute(
    "**ID__AFFARS_pgi_5301_topic_35__ID** ",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "**ID__AFFARS_AFMC_PGI_5301_9001_93__ID**", 2)

# Match the paragraph formatting used elsewhere in the document: a
# (borderless) paragraph border that just carries 5-twip spacing on
# all four sides, and a left indent of 225 twips (11.25pt) instead of
# the original 120 twips (6pt).
$p1 = $d.Paragraphs(1)
$p1.LeftIndent = 11.25

$b = $p1.Borders
$b.DistanceFromTop = 5
$b.DistanceFromBottom = 5
$b.DistanceFromLeft = 5
$b.DistanceFromRight = 5
